# Add team record (Wins / Losses / Ties) columns AD, AE, AF to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from the last existing header cell (AC1) onto
# the three new header cells so they match the rest of the header row
# (bold font, border, centered alignment).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every data row (2-49) gets the same team record for 2008 BOS.
$lastRow = 49
$wins = 95
$losses = 67
$ties = 0

$ws.Range("AD2:AD$lastRow").Value = $wins
$ws.Range("AE2:AE$lastRow").Value = $losses
$ws.Range("AF2:AF$lastRow").Value = $ties
